# Junction_Flooding_248.xlsx -- "custom accuracy + 데이터 1000개"
#
# 1) Round the 5th data row (columns B:AH) down to 2 decimal places
#    ("custom accuracy").
# 2) Drop the 6th data row entirely, shrinking the used range from
#    A1:AH6 down to A1:AH5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round a double to 2 decimal places (half-away-from-zero), working off the
# shortest round-trip decimal string rather than raw binary arithmetic so
# classic floating point noise (e.g. 9.0879999999999994) doesn't perturb the
# rounding decision.
function Round2([double]$value) {
    $isNegative = $false
    if ($value -lt 0) {
        $isNegative = $true
        $value = -$value
    }

    $text = $value.ToString()
    $dotIndex = $text.IndexOf(".")
    if ($dotIndex -lt 0) {
        $wholePart = $text
        $fractionPart = ""
    } else {
        $wholePart = $text.Substring(0, $dotIndex)
        $fractionPart = $text.Substring($dotIndex + 1)
    }

    # Make sure there are at least 3 fractional digits to inspect (the 3rd
    # one drives the rounding decision).
    $fractionPart = $fractionPart.PadRight(3, '0')
    $keepDigits = $fractionPart.Substring(0, 2)
    $roundingDigit = [int64]($fractionPart.Substring(2, 1))

    $hundredths = [int64]$wholePart * 100 + [int64]$keepDigits
    if ($roundingDigit -ge 5) {
        $hundredths = $hundredths + 1
    }

    $digits = ([string]$hundredths).PadLeft(3, '0')
    $len = $digits.Length
    $resultWhole = $digits.Substring(0, $len - 2)
    $resultFraction = $digits.Substring($len - 2)
    $result = [double]"$resultWhole.$resultFraction"

    if ($isNegative) {
        $result = -$result
    }
    return $result
}

# Columns B (2) through AH (34) on row 5.
$row = 5
$firstCol = 2
$lastCol = 34

for ($col = $firstCol; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item($row, $col)
    $original = [double]$cell.Value2
    $cell.Value = Round2($original)
}

# The sample only needs five rows of data now -- remove row 6 entirely
# (shifts the dimension from A1:AH6 to A1:AH5 automatically).
$ws.Rows("6:6").Delete()
